$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.035.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.09%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.835.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.05%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9978'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.32%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.59%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6223'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.83%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9979'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.30%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07491'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.93%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2935'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.61%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.34'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07692'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.85%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.834.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.008'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.21%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6746'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.41%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.65%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009365'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.46%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.956'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.06%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.030.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.082.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.05%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.79%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '221.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9981'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.23%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.158'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.30%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9981'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.32%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1400'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.520'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.36%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.86'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.74%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.491'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.65%  '

$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05642'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.16%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.174'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.09%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.138'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.05%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7471'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.19%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.843'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.44%  '

$ws.Range('E36').Value = '  -0.22%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.657'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.49%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.238.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.76%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.769'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01775'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.32%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.615'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8919'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.46%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9968'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.46%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.87%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.982.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.34%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.06%  '

$ws.Range('E47').Value = '  -1.00%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5074'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.95%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4064'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.67%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.059'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05838'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.52%  '
